$d = $word.ActiveDocument

# 1) Clear the RA number (keeps surrounding spaces) in the body text.
$d.Content.Find.Execute(" 000111746428 - 3 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "  ", 2)

# 2) Replace the placeholder tokens. MatchCase is required because
#    "QWER" / "QWR" / "Qwer" / "qwer" must be distinguished from one
#    another; ReplaceAll (2) handles the repeated occurrences in the
#    header.
$d.Content.Find.Execute("QWER", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TRE", 2)

$d.Content.Find.Execute("QWR", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TERE", 2)

$d.Content.Find.Execute("Qwer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Tre", 2)

$d.Content.Find.Execute("qwer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "tre", 2)
